$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.420.91'
$ws.Range('E2').Value = '  -4.28%  '
$ws.Range('D3').Value = '3.326.06'
$ws.Range('E3').Value = '  -4.71%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '547.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.610'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.21%  '
$ws.Range('D8').Value = '3.318.43'
$ws.Range('E8').Value = '  -4.83%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.611'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.152'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.05'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000265'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.86'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.94%  '
$ws.Range('D15').Value = '3.859.27'
$ws.Range('E15').Value = '  -5.36%  '
$ws.Range('D16').Value = '3.326.22'
$ws.Range('E16').Value = '  -5.46%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.116'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.91%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.74'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.22%  '
$ws.Range('D20').Value = '63.441.86'
$ws.Range('E20').Value = '  -4.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.967'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '405.40'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.23'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '82.72'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.52'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '28.97'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.36'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.27'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '569.92'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -9.07%  '
$ws.Range('E34').Value = '  -4.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.20'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.145'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.10'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.38'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0734'
$ws.Range('E40').Value = '  -9.12%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.150.92'
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.365'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0400'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.90%  '
$ws.Range('E48').Value = '  -5.09%  '
$ws.Range('E49').Value = '  -3.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.99%  '
